# Revert presentation to its earlier saved state ("Reverted to version 1"):
#   - remove the second slide ("Slie title") entirely
#   - clear the title text ("Sara Demo") on the remaining first slide

$p = $ppt.ActivePresentation

# Remove the second slide from the deck.
$p.Slides.Item(2).Delete()

# Clear the title placeholder's text on slide 1 (was "Sara Demo").
$s = $p.Slides.Item(1)
$title = $s.Shapes.Item(1)
$title.TextFrame2.TextRange.Paragraphs(1).Delete()
